$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price column cells being updated, to preserve
# literal formatting (leading/trailing zeros, dot-grouped values) as inline/shared strings
# rather than being reinterpreted as numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D41", "D43", "D45", "D48", "D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "48.204.19"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "2.509.82"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "109.21"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "320.57"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.545"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "39.88"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").Value = "20.08"
$ws.Range("E11").Value = "  +9.22%  "
$ws.Range("D12").Value = "0.0818"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "2.904.61"
$ws.Range("D16").Value = "2.520.86"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "48.056.25"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "0.0₃0943"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").Value = "2.73"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").Value = "274.05"
$ws.Range("E24").Value = "  +11.64%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "25.94"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  +4.88%  "
$ws.Range("D29").Value = "10.08"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.141"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "35.55"
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("D32").Value = "49.46"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "19.40"
$ws.Range("E33").Value = "  -5.88%  "
$ws.Range("D34").Value = "5.35"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "0.0784"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("D41").Value = "122.25"
$ws.Range("E41").Value = "  +4.20%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "21.72"
$ws.Range("E43").Value = "  -6.30%  "
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("D45").Value = "2.026.16"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").Value = "1.99"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("D51").Value = "79.45"
$ws.Range("E51").Value = "  +2.68%  "
